# Fixed a bug in weightBranch
# The rows of data (A2:F25) get reordered as a result of the fix.
# Apply the corrected values directly to the affected range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(902, 1, 0, 0, 0, 0),
    @(401, 9, 48, 67, 75, 45),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(901, 16, 15, 45, 60, 60),
    @(1203, 3, 15, 15, 15, 15),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(201, 9, 30, 15, 45, 30),
    @(1001, 18, 30, 75, 60, 72),
    @(601, 9, 60, 67, 60, 42),
    @(501, 9, 52, 30, 75, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(1, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2),
    @(1101, 0, 15, 30, 30, 0),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
